$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.246.20"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "2.382.42"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.37"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.83"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.92%  "
$ws.Range("D9").Value = "2.382.49"
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.87"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").Value = "2.837.32"
$ws.Range("E15").Value = "  +1.31%  "
$ws.Range("E16").Value = "  -1.25%  "
$ws.Range("D17").Value = "60.173.28"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").Value = "2.379.23"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.07"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +10.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.54"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.93"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.06"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.06"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("E25").Value = "  -2.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.26"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "562.28"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.49%  "
$ws.Range("E28").Value = "  -5.67%  "
$ws.Range("E29").Value = "  +0.36%  "
$ws.Range("D30").Value = "0.0₃0927"
$ws.Range("E30").Value = "  +1.65%  "
$ws.Range("E31").Value = "  +2.11%  "
$ws.Range("E32").Value = "  -2.12%  "
$ws.Range("E33").Value = "  -1.75%  "
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("E36").Value = "  +5.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.84"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.69%  "
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.09"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  +0.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.61"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.75%  "
$ws.Range("E45").Value = "  +4.29%  "
$ws.Range("D46").Value = "0.0₆0290"
$ws.Range("E46").Value = "  +3.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "140.55"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("E48").Value = "  +0.98%  "
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0502"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("E51").Value = "  -0.36%  "
